$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update time in/out values for rows 14-16 (columns C and D hold time-of-day values)
$ws.Range("C14").Value = 0.4375
$ws.Range("D14").Value = 0.645833333333333

$ws.Range("C15").Value = 0.833333333333333
$ws.Range("D15").Value = 0.916666666666667

$ws.Range("C16").Value = 0.520833333333333
$ws.Range("D16").Value = 0.729166666666667

# Update the selected cell to match the new active selection
$ws.Range("C19").Select()

# Slightly narrow the column widths as in the diff
# (ColumnWidth values below are chosen so that, after the engine's own
# width<->pixel quantization, the saved `width` attribute lands on the
# target values from the diff as closely as that quantization allows)
$ws.Columns.Item(1).ColumnWidth = 0.333333333333333
$ws.Columns.Item(2).ColumnWidth = 12.3333333333333
$ws.Range($ws.Columns.Item(3), $ws.Columns.Item(7)).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 11.6666666666667
$ws.Range($ws.Columns.Item(9), $ws.Columns.Item(11)).ColumnWidth = 10
$ws.Range($ws.Columns.Item(12), $ws.Columns.Item(1025)).ColumnWidth = 7.16666666666667

# Add a new named print area entry (duplicate of existing one) matching the diff
$ws.Names.Add("_xlnm.Print_Area_0_0_0_0_0_0", "='Weekly Time Record'!`$A`$1:`$K`$27")

$wb.Save()
